$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3102896511554718
$ws.Range("B1").Value = 0.2056666761636734
$ws.Range("C1").Value = 0.5889539122581482
$ws.Range("D1").Value = 3.634248495101929
$ws.Range("E1").Value = 3.843929529190063

$wb.Save()
